$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95, pushing existing rows 95-194 down to 96-195.
$ws.Rows("95:95").Insert()

# Populate the newly inserted row 95 with the new record's data.
$ws.Range("A95").Value = 7
$ws.Range("B95").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C95").Value = "Ñuble"
$ws.Range("D95").Value = 44664
$ws.Range("E95").Value = 16
$ws.Range("F95").Value = 100112017
$ws.Range("G95").Value = "Apio"
$ws.Range("H95").Value = "Americana (o)"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 100
$ws.Range("K95").Value = 8000
$ws.Range("L95").Value = 8500
$ws.Range("M95").Value = 8250
$ws.Range("N95").Value = "`$/docena de matas"
$ws.Range("O95").Value = "Provincia del Elquí"
$ws.Range("P95").Value = 1375
$ws.Range("Q95").Value = 6
$ws.Range("R95").Value = "Hortaliza"
